# Updated cryptos list on Fri May 26 18:46:21 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.903.30"
$ws.Range("E2").Value = "  +1.77%  "

# Row 3
$ws.Range("D3").Value = "1.844.81"
$ws.Range("E3").Value = "  +1.83%  "

# Row 4
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.38%  "

# Row 5
$ws.Range("D5").Value = "'309.39"
$ws.Range("E5").Value = "  +1.31%  "

# Row 6
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("D7").Value = "'0.4677"
$ws.Range("E7").Value = "  +3.84%  "

# Row 8
$ws.Range("D8").Value = "'0.3669"
$ws.Range("E8").Value = "  +2.38%  "

# Row 9
$ws.Range("D9").Value = "'0.07140"
$ws.Range("E9").Value = "  +1.17%  "

# Row 10
$ws.Range("D10").Value = "'0.9273"
$ws.Range("E10").Value = "  +4.31%  "

# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'19.56"
$ws.Range("E11").Value = "  +1.31%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07697"
$ws.Range("E12").Value = "  -1.00%  "

# Row 13
$ws.Range("D13").Value = "1.859.54"
$ws.Range("E13").Value = "  +3.06%  "

# Row 14
$ws.Range("D14").Value = "'5.284"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
$ws.Range("D15").Value = "'6.393"
$ws.Range("E15").Value = "  +1.41%  "

# Row 16
$ws.Range("D16").Value = "'88.06"
$ws.Range("E16").Value = "  +3.70%  "

# Row 17
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "  +0.32%  "

# Row 18
$ws.Range("D18").Value = "'0.000008626"
$ws.Range("E18").Value = "  +1.28%  "

# Row 19
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.45%  "

# Row 20
$ws.Range("D20").Value = "26.923.24"
$ws.Range("E20").Value = "  +1.65%  "

# Row 21
$ws.Range("D21").Value = "'14.39"
$ws.Range("E21").Value = "  +1.81%  "

# Row 22
$ws.Range("D22").Value = "'5.018"

# Row 23
$ws.Range("E23").Value = "  +1.59%  "

# Row 24
$ws.Range("D24").Value = "'1.933"
$ws.Range("E24").Value = "  -1.27%  "

# Row 25
$ws.Range("D25").Value = "'152.53"

# Row 26
$ws.Range("E26").Value = "  +2.82%  "

# Row 27
$ws.Range("D27").Value = "'2.030"
$ws.Range("E27").Value = "  -0.96%  "

# Row 28
$ws.Range("D28").Value = "'114.17"
$ws.Range("E28").Value = "  +1.65%  "

# Row 29
$ws.Range("D29").Value = "'4.889"
$ws.Range("E29").Value = "  +1.38%  "

# Row 30
$ws.Range("E30").Value = "  +2.12%  "

# Row 31
$ws.Range("D31").Value = "'3.208"
$ws.Range("E31").Value = "  +2.23%  "

# Row 32
$ws.Range("D32").Value = "'0.7475"

# Row 33
$ws.Range("D33").Value = "'1.172"
$ws.Range("E33").Value = "  +5.83%  "

# Row 34
$ws.Range("D34").Value = "'2.790"
$ws.Range("E34").Value = "  +2.05%  "

# Row 35
$ws.Range("D35").Value = "'4.461"
$ws.Range("E35").Value = "  +0.81%  "

# Row 36
$ws.Range("D36").Value = "'1.083"
$ws.Range("E36").Value = "  +1.59%  "

# Row 37
$ws.Range("D37").Value = "'0.01943"
$ws.Range("E37").Value = "  +0.91%  "

# Row 38
$ws.Range("D38").Value = "'2.968"
$ws.Range("E38").Value = "  +2.34%  "

# Row 39
$ws.Range("D39").Value = "'0.05190"
$ws.Range("E39").Value = "  +1.95%  "

# Row 40
$ws.Range("D40").Value = "'0.5212"
$ws.Range("E40").Value = "  +2.59%  "

# Row 41
$ws.Range("D41").Value = "'6.904"
$ws.Range("E41").Value = "  +2.37%  "

# Row 42
$ws.Range("D42").Value = "'0.1516"
$ws.Range("E42").Value = "  +0.68%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'8.134"
$ws.Range("E43").Value = "  +1.24%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'10.54"
$ws.Range("E44").Value = "  +5.41%  "

# Row 45
$ws.Range("D45").Value = "'0.4698"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("E46").Value = "  +0.52%  "

# Row 47
$ws.Range("D47").Value = "'100.66"
$ws.Range("E47").Value = "  +0.67%  "

# Row 48
$ws.Range("D48").Value = "'1.606"
$ws.Range("E48").Value = "  +2.07%  "

# Row 49
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'64.68"
$ws.Range("E49").Value = "  +1.89%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06028"
$ws.Range("E50").Value = "  +0.73%  "

# Row 51
$ws.Range("D51").Value = "'0.8870"
$ws.Range("E51").Value = "  +5.11%  "
